$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"10.93086033333333"
$ws.Cells.Item(2, 8).Value = [double]"32.792581"
$ws.Cells.Item(2, 9).Value = [double]"0.02611891973042464"
$ws.Cells.Item(2, 10).Value = [double]"0.02622656909968252"
$ws.Cells.Item(2, 13).Value = [double]"0.346153"
$ws.Cells.Item(2, 14).Value = [double]"1.038459"
$ws.Cells.Item(2, 15).Value = [double]"0.002587513800919532"
$ws.Cells.Item(2, 16).Value = [double]"0.002593123140481257"
$ws.Cells.Item(2, 17).Value = [double]"3.783750096964333"
$ws.Cells.Item(2, 18).Value = [double]"34.053750872679"
$ws.Cells.Item(2, 19).Value = [double]"6.758306526758323E-05"
$ws.Cells.Item(2, 20).Value = [double]"6.800872322781741E-05"
$ws.Cells.Item(3, 7).Value = [double]"10.93086033333333"
$ws.Cells.Item(3, 8).Value = [double]"32.792581"
$ws.Cells.Item(3, 9).Value = [double]"0.02611891973042464"
$ws.Cells.Item(3, 10).Value = [double]"0.02622656909968252"
$ws.Cells.Item(3, 15).Value = [double]"0.0001556606107424992"
$ws.Cells.Item(3, 16).Value = [double]"0.000155998059463248"
$ws.Cells.Item(3, 17).Value = [double]"0.2276242355813333"
$ws.Cells.Item(3, 18).Value = [double]"2.048618120232"
$ws.Cells.Item(3, 19).Value = [double]"4.065686997172213E-06"
$ws.Cells.Item(3, 20).Value = [double]"4.091293885929257E-06"
$ws.Cells.Item(4, 7).Value = [double]"10.93086033333333"
$ws.Cells.Item(4, 8).Value = [double]"32.792581"
$ws.Cells.Item(4, 9).Value = [double]"0.02611891973042464"
$ws.Cells.Item(4, 10).Value = [double]"0.02622656909968252"
$ws.Cells.Item(4, 13).Value = [double]"84.40796133333333"
$ws.Cells.Item(4, 14).Value = [double]"253.223884"
$ws.Cells.Item(4, 15).Value = [double]"0.630954418587972"
$ws.Cells.Item(4, 16).Value = [double]"0.6323222325801418"
$ws.Cells.Item(4, 17).Value = [double]"922.6516363560671"
$ws.Cells.Item(4, 18).Value = [double]"8303.864727204604"
$ws.Cells.Item(4, 19).Value = [double]"0.01647984781265599"
$ws.Cells.Item(4, 20).Value = [double]"0.01658364272602861"
$ws.Cells.Item(5, 7).Value = [double]"10.93086033333333"
$ws.Cells.Item(5, 8).Value = [double]"32.792581"
$ws.Cells.Item(5, 9).Value = [double]"0.02611891973042464"
$ws.Cells.Item(5, 10).Value = [double]"0.02622656909968252"
$ws.Cells.Item(5, 11).Value = [double]"2"
$ws.Cells.Item(5, 12).Value = [double]"1"
$ws.Cells.Item(5, 13).Value = [double]"0.868151"
$ws.Cells.Item(5, 14).Value = [double]"1.736302"
$ws.Cells.Item(5, 15).Value = [double]"0.006489479200764093"
$ws.Cells.Item(5, 16).Value = [double]"0.004335698275101748"
$ws.Cells.Item(5, 17).Value = [double]"9.489637329243667"
$ws.Cells.Item(5, 18).Value = [double]"56.937823975462"
$ws.Cells.Item(5, 19).Value = [double]"0.0001694981863370176"
$ws.Cells.Item(5, 20).Value = [double]"0.0001137104904073303"
$ws.Cells.Item(6, 7).Value = [double]"10.93086033333333"
$ws.Cells.Item(6, 8).Value = [double]"32.792581"
$ws.Cells.Item(6, 9).Value = [double]"0.02611891973042464"
$ws.Cells.Item(6, 10).Value = [double]"0.02622656909968252"
$ws.Cells.Item(6, 13).Value = [double]"48.13513433333333"
$ws.Cells.Item(6, 14).Value = [double]"144.405403"
$ws.Cells.Item(6, 15).Value = [double]"0.3598129277996019"
$ws.Cells.Item(6, 16).Value = [double]"0.3605929479448119"
$ws.Cells.Item(6, 17).Value = [double]"526.1584305239047"
$ws.Cells.Item(6, 18).Value = [double]"4735.425874715142"
$ws.Cells.Item(6, 19).Value = [double]"0.00939792497916688"
$ws.Cells.Item(6, 20).Value = [double]"0.009457115866132831"
$ws.Cells.Item(7, 9).Value = [double]"0.3422104954945279"
$ws.Cells.Item(7, 10).Value = [double]"0.3436209192170106"
$ws.Cells.Item(7, 13).Value = [double]"0.346153"
$ws.Cells.Item(7, 14).Value = [double]"1.038459"
$ws.Cells.Item(7, 15).Value = [double]"0.002587513800919532"
$ws.Cells.Item(7, 16).Value = [double]"0.002593123140481257"
$ws.Cells.Item(7, 17).Value = [double]"49.57475304774334"
$ws.Cells.Item(7, 18).Value = [double]"446.17277742969"
$ws.Cells.Item(7, 19).Value = [double]"0.0008854743799116024"
$ws.Cells.Item(7, 20).Value = [double]"0.0008910513571750707"
$ws.Cells.Item(8, 9).Value = [double]"0.3422104954945279"
$ws.Cells.Item(8, 10).Value = [double]"0.3436209192170106"
$ws.Cells.Item(8, 15).Value = [double]"0.0001556606107424992"
$ws.Cells.Item(8, 16).Value = [double]"0.000155998059463248"
$ws.Cells.Item(8, 19).Value = [double]"5.32686947311715E-05"
$ws.Cells.Item(8, 20).Value = [double]"5.360419658883116E-05"
$ws.Cells.Item(9, 9).Value = [double]"0.3422104954945279"
$ws.Cells.Item(9, 10).Value = [double]"0.3436209192170106"
$ws.Cells.Item(9, 13).Value = [double]"84.40796133333333"
$ws.Cells.Item(9, 14).Value = [double]"253.223884"
$ws.Cells.Item(9, 15).Value = [double]"0.630954418587972"
$ws.Cells.Item(9, 16).Value = [double]"0.6323222325801418"
$ws.Cells.Item(9, 17).Value = [double]"12088.59619406294"
$ws.Cells.Item(9, 18).Value = [double]"108797.3657465664"
$ws.Cells.Item(9, 19).Value = [double]"0.2159192242194517"
$ws.Cells.Item(9, 20).Value = [double]"0.2172791468005407"
$ws.Cells.Item(10, 9).Value = [double]"0.3422104954945279"
$ws.Cells.Item(10, 10).Value = [double]"0.3436209192170106"
$ws.Cells.Item(10, 11).Value = [double]"2"
$ws.Cells.Item(10, 12).Value = [double]"1"
$ws.Cells.Item(10, 13).Value = [double]"0.868151"
$ws.Cells.Item(10, 14).Value = [double]"1.736302"
$ws.Cells.Item(10, 15).Value = [double]"0.006489479200764093"
$ws.Cells.Item(10, 16).Value = [double]"0.004335698275101748"
$ws.Cells.Item(10, 17).Value = [double]"124.3333769551367"
$ws.Cells.Item(10, 18).Value = [double]"746.0002617308201"
$ws.Cells.Item(10, 19).Value = [double]"0.002220767892794913"
$ws.Cells.Item(10, 20).Value = [double]"0.00148983662673807"
$ws.Cells.Item(11, 9).Value = [double]"0.3422104954945279"
$ws.Cells.Item(11, 10).Value = [double]"0.3436209192170106"
$ws.Cells.Item(11, 13).Value = [double]"48.13513433333333"
$ws.Cells.Item(11, 14).Value = [double]"144.405403"
$ws.Cells.Item(11, 15).Value = [double]"0.3598129277996019"
$ws.Cells.Item(11, 16).Value = [double]"0.3605929479448119"
$ws.Cells.Item(11, 17).Value = [double]"6893.735999673414"
$ws.Cells.Item(11, 18).Value = [double]"62043.62399706072"
$ws.Cells.Item(11, 19).Value = [double]"0.1231317603076386"
$ws.Cells.Item(11, 20).Value = [double]"0.1239072802359679"
$ws.Cells.Item(12, 7).Value = [double]"157.1889546666667"
$ws.Cells.Item(12, 8).Value = [double]"471.566864"
$ws.Cells.Item(12, 9).Value = [double]"0.3755976715691904"
$ws.Cells.Item(12, 10).Value = [double]"0.3771457008466821"
$ws.Cells.Item(12, 13).Value = [double]"0.346153"
$ws.Cells.Item(12, 14).Value = [double]"1.038459"
$ws.Cells.Item(12, 15).Value = [double]"0.002587513800919532"
$ws.Cells.Item(12, 16).Value = [double]"0.002593123140481257"
$ws.Cells.Item(12, 17).Value = [double]"54.41142822473066"
$ws.Cells.Item(12, 18).Value = [double]"489.702854022576"
$ws.Cells.Item(12, 19).Value = [double]"0.0009718641587785219"
$ws.Cells.Item(12, 20).Value = [double]"0.0009779852441985527"
$ws.Cells.Item(13, 7).Value = [double]"157.1889546666667"
$ws.Cells.Item(13, 8).Value = [double]"471.566864"
$ws.Cells.Item(13, 9).Value = [double]"0.3755976715691904"
$ws.Cells.Item(13, 10).Value = [double]"0.3771457008466821"
$ws.Cells.Item(13, 15).Value = [double]"0.0001556606107424992"
$ws.Cells.Item(13, 16).Value = [double]"0.000155998059463248"
$ws.Cells.Item(13, 17).Value = [double]"3.273302791978666"
$ws.Cells.Item(13, 18).Value = [double]"29.459725127808"
$ws.Cells.Item(13, 19).Value = [double]"5.846576294992081E-05"
$ws.Cells.Item(13, 20).Value = [double]"5.883399746698906E-05"
$ws.Cells.Item(14, 7).Value = [double]"157.1889546666667"
$ws.Cells.Item(14, 8).Value = [double]"471.566864"
$ws.Cells.Item(14, 9).Value = [double]"0.3755976715691904"
$ws.Cells.Item(14, 10).Value = [double]"0.3771457008466821"
$ws.Cells.Item(14, 13).Value = [double]"84.40796133333333"
$ws.Cells.Item(14, 14).Value = [double]"253.223884"
$ws.Cells.Item(14, 15).Value = [double]"0.630954418587972"
$ws.Cells.Item(14, 16).Value = [double]"0.6323222325801418"
$ws.Cells.Item(14, 17).Value = [double]"13267.99920753109"
$ws.Cells.Item(14, 18).Value = [double]"119411.9928677798"
$ws.Cells.Item(14, 19).Value = [double]"0.2369850104879346"
$ws.Cells.Item(14, 20).Value = [double]"0.2384776115673763"
$ws.Cells.Item(15, 7).Value = [double]"157.1889546666667"
$ws.Cells.Item(15, 8).Value = [double]"471.566864"
$ws.Cells.Item(15, 9).Value = [double]"0.3755976715691904"
$ws.Cells.Item(15, 10).Value = [double]"0.3771457008466821"
$ws.Cells.Item(15, 11).Value = [double]"2"
$ws.Cells.Item(15, 12).Value = [double]"1"
$ws.Cells.Item(15, 13).Value = [double]"0.868151"
$ws.Cells.Item(15, 14).Value = [double]"1.736302"
$ws.Cells.Item(15, 15).Value = [double]"0.006489479200764093"
$ws.Cells.Item(15, 16).Value = [double]"0.004335698275101748"
$ws.Cells.Item(15, 17).Value = [double]"136.4637481828213"
$ws.Cells.Item(15, 18).Value = [double]"818.782489096928"
$ws.Cells.Item(15, 19).Value = [double]"0.002437433277503684"
$ws.Cells.Item(15, 20).Value = [double]"0.001635189964623"
$ws.Cells.Item(16, 7).Value = [double]"157.1889546666667"
$ws.Cells.Item(16, 8).Value = [double]"471.566864"
$ws.Cells.Item(16, 9).Value = [double]"0.3755976715691904"
$ws.Cells.Item(16, 10).Value = [double]"0.3771457008466821"
$ws.Cells.Item(16, 13).Value = [double]"48.13513433333333"
$ws.Cells.Item(16, 14).Value = [double]"144.405403"
$ws.Cells.Item(16, 15).Value = [double]"0.3598129277996019"
$ws.Cells.Item(16, 16).Value = [double]"0.3605929479448119"
$ws.Cells.Item(16, 17).Value = [double]"7566.311448596242"
$ws.Cells.Item(16, 18).Value = [double]"68096.80303736619"
$ws.Cells.Item(16, 19).Value = [double]"0.1351448978820237"
$ws.Cells.Item(16, 20).Value = [double]"0.1359960800730172"
$ws.Cells.Item(17, 7).Value = [double]"5.153359"
$ws.Cells.Item(17, 8).Value = [double]"10.306718"
$ws.Cells.Item(17, 9).Value = [double]"0.01231377640537609"
$ws.Cells.Item(17, 10).Value = [double]"0.008243018499152039"
$ws.Cells.Item(17, 13).Value = [double]"0.346153"
$ws.Cells.Item(17, 14).Value = [double]"1.038459"
$ws.Cells.Item(17, 15).Value = [double]"0.002587513800919532"
$ws.Cells.Item(17, 16).Value = [double]"0.002593123140481257"
$ws.Cells.Item(17, 17).Value = [double]"1.783850677927"
$ws.Cells.Item(17, 18).Value = [double]"10.703104067562"
$ws.Cells.Item(17, 19).Value = [double]"3.186206639034794E-05"
$ws.Cells.Item(17, 20).Value = [double]"2.137516201756623E-05"
$ws.Cells.Item(18, 7).Value = [double]"5.153359"
$ws.Cells.Item(18, 8).Value = [double]"10.306718"
$ws.Cells.Item(18, 9).Value = [double]"0.01231377640537609"
$ws.Cells.Item(18, 10).Value = [double]"0.008243018499152039"
$ws.Cells.Item(18, 15).Value = [double]"0.0001556606107424992"
$ws.Cells.Item(18, 16).Value = [double]"0.000155998059463248"
$ws.Cells.Item(18, 17).Value = [double]"0.107313547816"
$ws.Cells.Item(18, 18).Value = [double]"0.643881286896"
$ws.Cells.Item(18, 19).Value = [double]"1.916769955807419E-06"
$ws.Cells.Item(18, 20).Value = [double]"1.285894889987373E-06"
$ws.Cells.Item(19, 7).Value = [double]"5.153359"
$ws.Cells.Item(19, 8).Value = [double]"10.306718"
$ws.Cells.Item(19, 9).Value = [double]"0.01231377640537609"
$ws.Cells.Item(19, 10).Value = [double]"0.008243018499152039"
$ws.Cells.Item(19, 13).Value = [double]"84.40796133333333"
$ws.Cells.Item(19, 14).Value = [double]"253.223884"
$ws.Cells.Item(19, 15).Value = [double]"0.630954418587972"
$ws.Cells.Item(19, 16).Value = [double]"0.6323222325801418"
$ws.Cells.Item(19, 17).Value = [double]"434.9845272087853"
$ws.Cells.Item(19, 18).Value = [double]"2609.907163252712"
$ws.Cells.Item(19, 19).Value = [double]"0.007769431632476357"
$ws.Cells.Item(19, 20).Value = [double]"0.005212243860583226"
$ws.Cells.Item(20, 7).Value = [double]"5.153359"
$ws.Cells.Item(20, 8).Value = [double]"10.306718"
$ws.Cells.Item(20, 9).Value = [double]"0.01231377640537609"
$ws.Cells.Item(20, 10).Value = [double]"0.008243018499152039"
$ws.Cells.Item(20, 11).Value = [double]"2"
$ws.Cells.Item(20, 12).Value = [double]"1"
$ws.Cells.Item(20, 13).Value = [double]"0.868151"
$ws.Cells.Item(20, 14).Value = [double]"1.736302"
$ws.Cells.Item(20, 15).Value = [double]"0.006489479200764093"
$ws.Cells.Item(20, 16).Value = [double]"0.004335698275101748"
$ws.Cells.Item(20, 17).Value = [double]"4.473893769209"
$ws.Cells.Item(20, 18).Value = [double]"17.895575076836"
$ws.Cells.Item(20, 19).Value = [double]"7.990999586554776E-05"
$ws.Cells.Item(20, 20).Value = [double]"3.57392410884053E-05"
$ws.Cells.Item(21, 7).Value = [double]"5.153359"
$ws.Cells.Item(21, 8).Value = [double]"10.306718"
$ws.Cells.Item(21, 9).Value = [double]"0.01231377640537609"
$ws.Cells.Item(21, 10).Value = [double]"0.008243018499152039"
$ws.Cells.Item(21, 13).Value = [double]"48.13513433333333"
$ws.Cells.Item(21, 14).Value = [double]"144.405403"
$ws.Cells.Item(21, 15).Value = [double]"0.3598129277996019"
$ws.Cells.Item(21, 16).Value = [double]"0.3605929479448119"
$ws.Cells.Item(21, 17).Value = [double]"248.0576277328923"
$ws.Cells.Item(21, 18).Value = [double]"1488.345766397354"
$ws.Cells.Item(21, 19).Value = [double]"0.004430655940688028"
$ws.Cells.Item(21, 20).Value = [double]"0.002972374340572852"
$ws.Cells.Item(22, 7).Value = [double]"102.0140613333333"
$ws.Cells.Item(22, 8).Value = [double]"306.042184"
$ws.Cells.Item(22, 9).Value = [double]"0.2437591368004809"
$ws.Cells.Item(22, 10).Value = [double]"0.2447637923374727"
$ws.Cells.Item(22, 13).Value = [double]"0.346153"
$ws.Cells.Item(22, 14).Value = [double]"1.038459"
$ws.Cells.Item(22, 15).Value = [double]"0.002587513800919532"
$ws.Cells.Item(22, 16).Value = [double]"0.002593123140481257"
$ws.Cells.Item(22, 17).Value = [double]"35.31247337271733"
$ws.Cells.Item(22, 18).Value = [double]"317.812260354456"
$ws.Cells.Item(22, 19).Value = [double]"0.0006307301305714764"
$ws.Cells.Item(22, 20).Value = [double]"0.0006347026538622493"
$ws.Cells.Item(23, 7).Value = [double]"102.0140613333333"
$ws.Cells.Item(23, 8).Value = [double]"306.042184"
$ws.Cells.Item(23, 9).Value = [double]"0.2437591368004809"
$ws.Cells.Item(23, 10).Value = [double]"0.2447637923374727"
$ws.Cells.Item(23, 15).Value = [double]"0.0001556606107424992"
$ws.Cells.Item(23, 16).Value = [double]"0.000155998059463248"
$ws.Cells.Item(23, 17).Value = [double]"2.124340813205333"
$ws.Cells.Item(23, 18).Value = [double]"19.119067318848"
$ws.Cells.Item(23, 19).Value = [double]"3.794369610842727E-05"
$ws.Cells.Item(23, 20).Value = [double]"3.818267663151115E-05"
$ws.Cells.Item(24, 7).Value = [double]"102.0140613333333"
$ws.Cells.Item(24, 8).Value = [double]"306.042184"
$ws.Cells.Item(24, 9).Value = [double]"0.2437591368004809"
$ws.Cells.Item(24, 10).Value = [double]"0.2447637923374727"
$ws.Cells.Item(24, 13).Value = [double]"84.40796133333333"
$ws.Cells.Item(24, 14).Value = [double]"253.223884"
$ws.Cells.Item(24, 15).Value = [double]"0.630954418587972"
$ws.Cells.Item(24, 16).Value = [double]"0.6323222325801418"
$ws.Cells.Item(24, 17).Value = [double]"8610.798944480295"
$ws.Cells.Item(24, 18).Value = [double]"77497.19050032266"
$ws.Cells.Item(24, 19).Value = [double]"0.1538009044354533"
$ws.Cells.Item(24, 20).Value = [double]"0.1547695876256129"
$ws.Cells.Item(25, 7).Value = [double]"102.0140613333333"
$ws.Cells.Item(25, 8).Value = [double]"306.042184"
$ws.Cells.Item(25, 9).Value = [double]"0.2437591368004809"
$ws.Cells.Item(25, 10).Value = [double]"0.2447637923374727"
$ws.Cells.Item(25, 11).Value = [double]"2"
$ws.Cells.Item(25, 12).Value = [double]"1"
$ws.Cells.Item(25, 13).Value = [double]"0.868151"
$ws.Cells.Item(25, 14).Value = [double]"1.736302"
$ws.Cells.Item(25, 15).Value = [double]"0.006489479200764093"
$ws.Cells.Item(25, 16).Value = [double]"0.004335698275101748"
$ws.Cells.Item(25, 17).Value = [double]"88.56360936059468"
$ws.Cells.Item(25, 18).Value = [double]"531.3816561635681"
$ws.Cells.Item(25, 19).Value = [double]"0.00158186984826293"
$ws.Cells.Item(25, 20).Value = [double]"0.001061221952244943"
$ws.Cells.Item(26, 7).Value = [double]"102.0140613333333"
$ws.Cells.Item(26, 8).Value = [double]"306.042184"
$ws.Cells.Item(26, 9).Value = [double]"0.2437591368004809"
$ws.Cells.Item(26, 10).Value = [double]"0.2447637923374727"
$ws.Cells.Item(26, 13).Value = [double]"48.13513433333333"
$ws.Cells.Item(26, 14).Value = [double]"144.405403"
$ws.Cells.Item(26, 15).Value = [double]"0.3598129277996019"
$ws.Cells.Item(26, 16).Value = [double]"0.3605929479448119"
$ws.Cells.Item(26, 17).Value = [double]"4910.460546168905"
$ws.Cells.Item(26, 18).Value = [double]"44194.14491552015"
$ws.Cells.Item(26, 19).Value = [double]"0.0877076886900847"
$ws.Cells.Item(26, 20).Value = [double]"0.08826009742912104"
